$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New "Testomgång1" label in the previously-empty row 1 ---
$ws.Range("E1").Value = "Testomgång1"
$ws.Range("E1").Font.Bold = $true

# --- Make room for the new testfall 13 row (testfall 12 fits into the
#     already-empty row 14 gap, so only one row needs inserting, right
#     before the old row 15 "Slutsats" block) ---
$ws.Rows("15:15").Insert()

# --- Testfall 12: "Likbent triangel där lika sidor kommer sist som inparametrar" ---
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Likbent triangel där lika sidor kommer sist som inparametrar"
$ws.Range("C14").Value = "4,2 3 3"
$ws.Range("D14").Value = "Triangeln är likbent"
$ws.Range("E14").Value = "Pass"

# --- Testfall 13: "En inparameter ej giltig/ej double" ---
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "En inparameter ej giltig/ej double"
$ws.Range("C15").Value = "4,2 3 ett"
$ws.Range("D15").Value = "Hanterat fel"
$ws.Range("E15").Value = "Fail"
$ws.Range("F15").Value = "Unhandled Exception: System.FormatException: Input string was not in a correct format. - Ohanterat fel, applikationen hänger sig"
$ws.Range("F15").WrapText = $true
$ws.Rows("15:15").RowHeight = 30

# --- Selection matches the author's final cursor position ---
[void]$ws.Range("F15").Select()
